$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (column C) date from 45443 (2024-05-31) to 45444 (2024-06-01)
# for every data row that stays in the sheet (rows 2 through 28).
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = 45444
}

# Row 28 loses its explicit custom row height (ht="15" customHeight="1")
# and reverts to the sheet's default row height.
$ws.Rows.Item(28).AutoFit()

# Row 29 (A 21829-2024) is removed entirely.
$ws.Rows.Item(29).Delete()
